# Adds new categorized links to the "Arkusz1" sheet (rows 506-521),
# and applies a distinct "Courier New" style to A507 ("dramatic mode").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 506; Link = "https://medium.com/content-uneditable/circular-dependencies-in-javascript-a-k-a-coding-is-not-a-rock-paper-scissors-game-9c2a9eccd4bc"; Category = "IT" },
    @{ Row = 507; Link = "https://medium.com/javascript-scene/inside-the-dev-team-death-spiral-6a7ea255467b"; Category = "IT" },
    @{ Row = 508; Link = "https://qz.com/1285418/giant-predatory-worms-from-asia-are-invading-france/?utm_source=parAO"; Category = "Nature" },
    @{ Row = 509; Link = "https://sexier.com/live-sex-chats/?queryid=138"; Category = "Porn" },
    @{ Row = 510; Link = "https://vividcams.com/?AFNO=MjAyOTEzLjU0LjQ2LjEwMi4xOS4wLjAuMC4w&mobile=0&nats=MjAyOTEzLjU0LjQ2LjEwMi4xOS4wLjAuMC4w&strack=0&switched=1"; Category = "Porn" },
    @{ Row = 511; Link = "https://cams.com/go/g1424946-pct"; Category = "Porn" },
    @{ Row = 512; Link = "http://myporncams.com/"; Category = "Porn" },
    @{ Row = 513; Link = "https://www.xlovecam.com/en/"; Category = "Porn" },
    @{ Row = 514; Link = "https://livecam-experts.com/"; Category = "Porn" },
    @{ Row = 515; Link = "https://www.flirt4free.com/live/girls/"; Category = "Porn" },
    @{ Row = 516; Link = "https://www.pornication.com/?AFNO=1-247331-2-2-bestxxxsites"; Category = "Porn" },
    @{ Row = 517; Link = "https://www.myfreecams.com/?cam=30546&omp=2&track=102c530982c0f2496e5886e66ee166&skip_oapopup=1&r=0&mfwd=1#{model}"; Category = "Porn" },
    @{ Row = 518; Link = "https://www.privatefeeds.com/?AFNO=1-247331-2-bestxxxsites"; Category = "Porn" },
    @{ Row = 519; Link = "https://www.streamate.com/?AFNO=1-0-642160-356079&DF=0&UHNSMTY=303"; Category = "Porn" },
    @{ Row = 520; Link = "https://www.evilangellive.com/?AFNO=1-247331-2-2-bestxxxsites"; Category = "Porn" },
    @{ Row = 521; Link = "https://www.watchmygf.me/"; Category = "Porn" }
)

foreach ($item in $rows) {
    $ws.Cells.Item($item.Row, 1).Value = $item.Link
    $ws.Cells.Item($item.Row, 2).Value = $item.Category
}

# Make row 507 ("the ER" entry) errorproof / visually distinct: dramatic,
# monospaced "Courier New" styling matching the rest of the quote-prefixed
# vertical-centered cells already used in the sheet (vertical-center +
# black text), just swapped to a fixed-width font.
$target = $ws.Range("A507")
$ws.Range("A488").Copy()
$target.PasteSpecial(-4122)
$target.Font.Size = 10
$target.Font.Name = "Courier New"

$ws.Range("A521").Select()
